# jangseong DataImport.xlsx correction
# Re-imports the roster data (old 3-row dataset replaced with corrected
# 5-row dataset), applies center alignment / text format to the "class code"
# column, switches the sheet font to Malgun Gothic (맑은 고딕), and moves the
# active selection below the refreshed data, matching the
# "Old files delete, jangseong sql file correction" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Extend the table to 5 rows, re-using row 3's existing formatting so
#    the new rows start out identically styled to the rest of the data.
# ---------------------------------------------------------------------
$ws.Range("A3:F3").Copy($ws.Range("A4:F4"))
$ws.Range("A3:F3").Copy($ws.Range("A5:F5"))

# ---------------------------------------------------------------------
# 2. Write the corrected data grid (5 rows x 6 columns: No / Class code /
#    Region / Date / Student No / Name)
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = 1
$ws.Cells.Item(1, 2).Value = "23-01"
$ws.Cells.Item(1, 3).Value = "울산"
$ws.Cells.Item(1, 4).Value = 45213
$ws.Cells.Item(1, 5).Value = "23-000001"
$ws.Cells.Item(1, 6).Value = "김채원"

$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "23-01"
$ws.Cells.Item(2, 3).Value = "울산"
$ws.Cells.Item(2, 4).Value = 45213
$ws.Cells.Item(2, 5).Value = "23-000002"
$ws.Cells.Item(2, 6).Value = "허윤진"

$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "23-01"
$ws.Cells.Item(3, 3).Value = "울산"
$ws.Cells.Item(3, 4).Value = 44999
$ws.Cells.Item(3, 5).Value = "23-000003"
$ws.Cells.Item(3, 6).Value = "홍은채"

$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(4, 2).Value = "23-02"
$ws.Cells.Item(4, 3).Value = "울산"
$ws.Cells.Item(4, 4).Value = 45000
$ws.Cells.Item(4, 5).Value = "23-000004"
$ws.Cells.Item(4, 6).Value = "카즈하"

$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "23-03"
$ws.Cells.Item(5, 3).Value = "울산"
$ws.Cells.Item(5, 4).Value = 45001
$ws.Cells.Item(5, 5).Value = "23-000005"
$ws.Cells.Item(5, 6).Value = "사쿠라"

# ---------------------------------------------------------------------
# 3. Formatting across the refreshed range (A1:F5)
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:F5")

# Whole block becomes horizontally (in addition to vertically) centered
$dataRange.HorizontalAlignment = -4108   # xlCenter
$dataRange.VerticalAlignment = -4108     # xlCenter

# Class-code column (B) is stored as text so codes like "23-01" keep
# their leading digits / dashes intact
$ws.Range("B1:B5").NumberFormat = "@"

# Sheet font switched from Calibri to Malgun Gothic
$dataRange.Font.Name = "맑은 고딕"

# ---------------------------------------------------------------------
# 4. Move the active selection to just below the refreshed data
# ---------------------------------------------------------------------
$ws.Range("F6").Select()
